$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 60; existing rows 60-83 shift down to 61-84.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly price record.
$ws.Cells.Item(60, 1).Value = 4
$ws.Cells.Item(60, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(60, 3).Value = "Los Lagos"
$ws.Cells.Item(60, 4).Value = 45009
$ws.Cells.Item(60, 5).Value = 10
$ws.Cells.Item(60, 6).Value = 100112043
$ws.Cells.Item(60, 7).Value = "Pepino dulce"
$ws.Cells.Item(60, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 100
$ws.Cells.Item(60, 11).Value = 17000
$ws.Cells.Item(60, 12).Value = 18000
$ws.Cells.Item(60, 13).Value = 17500
$ws.Cells.Item(60, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(60, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value = 972
$ws.Cells.Item(60, 17).Value = 18
$ws.Cells.Item(60, 18).Value = "Hortaliza"
